$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text cells whose new value looks like a plain number need to be force-written
# as text (NumberFormat "@" while assigning, then restored to General/Normal style)
# so they keep matching the original inlineStr (no numeric coercion, no residual style).
$textCells = @{
    "D5" = "294.07"
    "D6" = "87.71"
    "D10" = "30.83"
    "D11" = "50.99"
    "D12" = "0.0783"
    "D14" = "6.46"
    "D16" = "13.85"
    "D21" = "11.29"
    "D23" = "65.77"
    "D24" = "236.66"
    "D25" = "1.00"
    "D28" = "23.36"
    "D29" = "9.30"
    "D31" = "157.74"
    "D32" = "31.83"
    "D35" = "3.04"
    "D36" = "0.0717"
    "D37" = "2.32"
    "D40" = "0.0992"
    "D41" = "15.47"
    "D44" = "18.45"
    "D45" = "10.15"
    "D48" = "2.73"
}

foreach ($cellRef in $textCells.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $textCells[$cellRef]
    $rng.Style = "Normal"
}

# Remaining cells are already non-numeric-looking strings; plain assignment is safe.
$ws.Range("D2").Value = "40.176.92"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "2.227.42"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "2.574.50"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "2.241.68"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "40.114.19"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  -6.65%  "
$ws.Range("E31").Value = "  +3.02%  "
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("E35").Value = "  +7.13%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").Value = "2.086.74"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  +5.24%  "
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  -10.84%  "
$ws.Range("E48").Value = "  +3.00%  "
$ws.Range("D49").Value = "2.447.08"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("E51").Value = "  +3.84%  "
